$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new rows (line7, line8) are being added to the table, pushing all the
# "extr*" rows down by two positions. Rows 2-7 (line1..line6) stay untouched.
# Build the full target table for rows 8-17 (values already reflect the
# post-insert state) and write it directly, cell by cell.

# First, extend formatting for the two brand-new rows (16 and 17) by copying
# the style of the last existing data row (15) which carries the bold /
# bordered / centered style used for column A.
$ws.Range("A15:E15").Copy()
$ws.Range("A16:E16").PasteSpecial(-4122)
$ws.Range("A15:E15").Copy()
$ws.Range("A17:E17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$rows = @(
  @{ Row=8;  A=6;  B="line7"; C=14; D=11; E=$true  },
  @{ Row=9;  A=7;  B="line8"; C=16; D=9;  E=$true  },
  @{ Row=10; A=8;  B="extr1"; C=5;  D=12; E=$true  },
  @{ Row=11; A=9;  B="extr2"; C=5;  D=9;  E=$true  },
  @{ Row=12; A=10; B="extr3"; C=10; D=11; E=$false },
  @{ Row=13; A=11; B="extr4"; C=7;  D=8;  E=$true  },
  @{ Row=14; A=12; B="extr5"; C=9;  D=11; E=$true  },
  @{ Row=15; A=13; B="extr6"; C=7;  D=11; E=$true  },
  @{ Row=16; A=14; B="extr7"; C=5;  D=7;  E=$true  },
  @{ Row=17; A=15; B="extr8"; C=8;  D=5;  E=$false }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.A
    $ws.Range("B$n").Value = $r.B
    $ws.Range("C$n").Value = $r.C
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
}
